# Refactor security vulnerability checks: add a new validation-record row
# (row 57) to each of the four MID_* log sheets, mirroring the existing
# row 56 pattern/format but with updated input-validation data.

$wb = $excel.ActiveWorkbook

# New row-57 data per sheet (A..I), in sheet order matching the workbook.
$newRows = @(
    @{
        A = 45843.46112268518
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x64"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 356
        I = 7
    },
    @{
        A = 45843.46112268518
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x5C"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 348
        I = 25
    },
    @{
        A = 45843.46112268518
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x68"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 104
        I = 15
    },
    @{
        A = 45843.46112268518
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x7D"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 125
        I = 9
    }
)

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $row = 57
    $data = $newRows[$i - 1]

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $cellA.Value = $data.A

    $ws.Cells.Item($row, 2).Value = $data.B
    $ws.Cells.Item($row, 3).Value = $data.C
    $ws.Cells.Item($row, 4).Value = $data.D
    $ws.Cells.Item($row, 5).Value = $data.E
    $ws.Cells.Item($row, 6).Value = $data.F
    $ws.Cells.Item($row, 7).Value = $data.G
    $ws.Cells.Item($row, 8).Value = $data.H
    $ws.Cells.Item($row, 9).Value = $data.I
}
